$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.150.55"
$ws.Range("E2").Value = "  -1.48%  "

# Row 3
$ws.Range("D3").Value = "1.782.62"
$ws.Range("E3").Value = "  -1.74%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.42%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.02"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.99%  "

# Row 6
$ws.Range("E6").Value = "  +0.34%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3828"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.26%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3424"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.60%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.13"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.46%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.193"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.18%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07491"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.62%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.34%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.73"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.10%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.443"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.62%  "

# Row 15
$ws.Range("D15").Value = "1.784.91"
$ws.Range("E15").Value = "  -1.59%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.075"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.92%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001091"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.40%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06682"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.62"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.41%  "

# Row 20
$ws.Range("E20").Value = "  +0.31%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.614"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.34"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.41%  "

# Row 23
$ws.Range("D23").Value = "27.156.40"
$ws.Range("E23").Value = "  -1.41%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -6.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.384"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.39%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.531"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.63%  "

# Row 27
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.466"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.10%  "

# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.21"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.20%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.10"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.11%  "

# Row 30
$ws.Range("D30").Value = "1.990.09"
$ws.Range("E30").Value = "  -1.36%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "134.07"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.90%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.020"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.61%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.048"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08710"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.36%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.23"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.42%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.648"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.28%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6881"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.31%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.406"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.25%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2197"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.24%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06313"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.63%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.758"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.93%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02331"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.61%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.236"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.92%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.38"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.38%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6468"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.74%  "

# Row 46
$ws.Range("E46").Value = "  +0.36%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.858"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.38%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.140"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.04%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.21"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.12%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07127"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.22%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.98"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.24%  "
